$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rows 12-18: remove the "US " / "US" notes in column A, add an order date in column J ---
$ws.Range("A12").ClearContents()
$ws.Range("A13").ClearContents()
$ws.Range("A14").ClearContents()
$ws.Range("A15").ClearContents()
$ws.Range("A16").ClearContents()
$ws.Range("A17").ClearContents()
$ws.Range("A18").ClearContents()

$ws.Range("J12").Value = 41684
$ws.Range("J12").NumberFormat = "m/d/yy"
$ws.Range("J13").Value = 41684
$ws.Range("J13").NumberFormat = "m/d/yy"
$ws.Range("J14").Value = 41684
$ws.Range("J14").NumberFormat = "m/d/yy"
$ws.Range("J15").Value = 41684
$ws.Range("J15").NumberFormat = "m/d/yy"
$ws.Range("J16").Value = 41684
$ws.Range("J16").NumberFormat = "m/d/yy"
$ws.Range("J17").Value = 41684
$ws.Range("J17").NumberFormat = "m/d/yy"
$ws.Range("J18").Value = 41684
$ws.Range("J18").NumberFormat = "m/d/yy"

# --- Row 19: motor changed from Turnigy 450 H2218 to NTM Prop Drive; remove "International" note, add order date ---
$ws.Range("A19").ClearContents()
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 18.48
$ws.Range("G19").Value = "http://hobbyking.com/hobbyking/store/__26486__NTM_Prop_Drive_Series_35_36A_1800Kv_875w_US_Warehouse_.html"
$ws.Range("J19").Value = 41684
$ws.Range("J19").NumberFormat = "m/d/yy"

# --- Row 37: returned the remote control ---
$ws.Range("D37").Value = 0

# --- Row 38: new ESC - Speed Controllers part replacing the placeholder "ESC" row ---
$ws.Range("G38").Value = "http://www.hobbyexpress.com/erc_rapid_drive_35a_brushless_esc_1039744_prd1.htm"
$ws.Range("C38").Value = "35A"
$ws.Range("J37").Value = "Returned"
$ws.Range("B38").Value = "ESC - Speed Controllers"
$ws.Range("D38").Value = 5
$ws.Range("E38").Value = 29.99

# --- Update the view to match where the editor left off ---
$ws.Range("B39").Select()
$excel.ActiveWindow.ScrollRow = 7
